# Update the "Metadata" sheet: URL, Version, Date, Publisher
# (moving from the "ibm.com / Alvearie" project identity to "linuxforhealth.org / LinuxForHealth")
$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/communication-vendor"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# Update the "Elements" sheet:
#  - Extension.url's Fixed Value should track the same URL rename as above
#  - The mis-duplicated ele-1/ext-1 constraint on the root "Extension" row is removed
#    (that constraint correctly lives only on the "Extension.extension" row)
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/communication-vendor"
$elements.Range("AI2").Value = ""
